# Comparison momentum error and energy error for the inviscid shear layer test case.
# Restructures the "testcase_list" sheet: inserts an "Initial condition" column,
# adds two new ROM rows, an "Other" column, and wraps the range in an Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$e100 = [double]"1E+100"

# --- 1. Insert a new column at D ("Initial condition"), shifting old D..J to E..K ---
$ws.Columns.Item(4).Insert()

# --- 2. Row 4 (was the 200x200 FOM case) becomes ID 1 / shear_layer01 ---
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "shear_layer01"
$ws.Range("D3").Value = "Initial condition"
$ws.Range("D4").Value = "shear layer + 1 for u comp"
$ws.Range("F4").Value = "40x40"

# --- 3. Row 5 (was the 40x40 momcons FOM case) becomes ID 2 / shear_layer02 ---
$ws.Range("A5").Value = 2
$ws.Range("D5").Value = "shear layer standard"
$ws.Range("B5").Value = "shear_layer02"
$ws.Range("F5").Value = "200x200"
$ws.Range("I5").Value = 6

# --- 4. New row 6: shear_layer03 (copy row 4 as a starting template, then adjust) ---
$ws.Range("A4:K4").Copy()
$ws.Range("A6").PasteSpecial(-4163)
$ws.Range("E6").NumberFormat = "0.00E+00"
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "shear_layer03"
$ws.Range("F6").Value = "200x200"

# --- 5. New row 7: ROM case, not momentum conserving (copy row 4 as template) ---
$ws.Range("A4:K4").Copy()
$ws.Range("A7").PasteSpecial(-4163)
$ws.Range("E7").NumberFormat = "0.00E+00"
$ws.Range("A7").Value = 4
$ws.Range("B7").ClearContents()
$ws.Range("J7").Value = "GL1"
$ws.Range("K7").Value = "ROM"
$ws.Range("M7").Value = "shear_layer01/matlab_data.mat"
$ws.Range("N3").Value = "Other"
$ws.Range("N7").Value = "not momentum conserving"
$ws.Range("L7").Value = "2,4,8,16"

# --- 6. New row 8: ROM case, momentum conserving (copy row 7 as template) ---
$ws.Range("A7:N7").Copy()
$ws.Range("A8").PasteSpecial(-4163)
$ws.Range("E8").NumberFormat = "0.00E+00"
$ws.Range("A8").Value = 5
$ws.Range("N8").Value = "momentum conserving"

$excel.CutCopyMode = 0

# --- 7. Column widths ---
$ws.Columns.Item(2).ColumnWidth = 37.666666666666664   # B Name
$ws.Columns.Item(4).ColumnWidth = 24.666666666666668   # D Initial condition (same as C)
$ws.Columns.Item(5).ColumnWidth = 17.0                 # E Reynolds number
$ws.Columns.Item(7).ColumnWidth = 19.166666666666668   # G boundary conditions
$ws.Columns.Item(10).ColumnWidth = 22.833333333333332  # J time integration method
$ws.Columns.Item(13).ColumnWidth = 29.166666666666668  # M snapshot matrix file
$ws.Columns.Item(14).ColumnWidth = 29.166666666666668  # N Other

# --- 8. Wrap the range in a table ---
$rng = $ws.Range("A3:N8")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# --- 9. Selection / view ---
$ws.Range("N9").Select()

$wb.Save()
